# A02QVZ2_Fin Buff Calc.xlsx - "Add files via upload" revision
#
# Updates the three input cells on Sheet1 (Gross Expenditures, Total M,
# Total Labor Cost) with the new figures from revision 4 of the workbook,
# then moves the active selection to D5 the way the author last left it
# before saving (the dependent formulas in D6:D9/E8:E9 recalculate
# automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = 405431.51
$ws.Range("D4").Value = 745.73
$ws.Range("D5").Value = 90657.57

$ws.Range("D5").Select()
